$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: Predicted (Robotic Arthroplasty query) ---
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A6").Value = "Predicted"

$b6 = "`n" + '"digital surgery" OR "surgical navigation" OR "haptic feedback" OR "robot-assisted surgery" OR "prosthetic devices" OR "computer-assisted surgery" OR "surgical workflow" OR (robotic arthroplasty) OR ' + "`n" + '(("clinical trials" OR "rehabilitation" OR "recovery time" OR "3d imaging" OR "patient satisfaction" OR "minimally invasive surgery" OR "surgical safety" OR "advanced imaging techniques" OR "orthopedic surgery" OR "total hip arthroplasty" OR "total knee arthroplasty" OR "biomechanics" OR "surgeon training" OR "motion planning" OR "joint replacement" OR "surgical precision" OR "robotic surgery" OR "implant technology" OR "surgery simulation" OR "surgical robotics") AND ("Arthoplasty"))' + "`n"
$ws.Range("B6").Value = $b6
$ws.Range("C6").Value = 0.957
$ws.Range("D6").Value = 0.061
$ws.Range("E6").Value = 0.244
$ws.Range("F6").Value = 0.33
$ws.Range("G6").Value = 0.694
$ws.Rows.Item(6).AutoFit()

# --- Row 7: Baseline (Robotic Arthroplasty) ---
$ws.Range("A3").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A7").Value = "Baseline"

$ws.Range("B7").Value = "Robotic Arthroplasty"
$ws.Range("C7").Value = 0.957
$ws.Range("D7").Value = 0.595
$ws.Range("E7").Value = 0.853
$ws.Range("F7").Value = 0.401
$ws.Range("G7").Value = 0.749

# --- Row 8: Predicted (Soft Robotics query) ---
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Predicted"

$b8 = "`n" + '"compliant materials" OR (soft robotics) OR "soft actuators" OR "rehabilitation robotics" OR "pneumatic actuators" OR "robotic locomotion" OR "soft robotic systems" OR "soft robot control" OR "multi-material printing" OR "bio-inspired robotics" OR "continuum robots" OR "soft robot fabrication" OR "shape-morphing structures" OR "energy-efficient robotics" OR "soft robotic grippers" OR "soft exoskeletons" OR "soft robotic arms" OR "soft wearable robots" OR "autonomous soft robots" OR ' + "`n" + '(("material properties" OR "adaptive control" OR "sensor integration" OR "flexible materials" OR "human-robot interaction" OR "deformation mechanics" OR "lightweight structures" OR "soft sensors" OR "artificial muscles" OR "robotic manipulation") AND (Robot OR Soft))' + "`n"
$ws.Range("B8").Value = $b8
$ws.Range("C8").Value = 0.722
$ws.Range("D8").Value = 0.212
$ws.Range("E8").Value = 0.487
$ws.Range("F8").Value = 0.451
$ws.Range("G8").Value = 0.644
$ws.Rows.Item(8).AutoFit()

# --- Row 9: Baseline (Soft Robotics) ---
$ws.Range("A3").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "Baseline"

$ws.Range("B9").Value = "Soft Robotics"
$ws.Range("C9").Value = 0.556
$ws.Range("D9").Value = 0.452
$ws.Range("E9").Value = 0.531
$ws.Range("F9").Value = 0.492
$ws.Range("G9").Value = 0.542

# --- Row 10: Predicted (Crop Yield Prediction query) ---
$ws.Range("A2").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Predicted"

$b10 = "`n" + '"crop yield estimation" OR "crop simulation models" OR "yield forecasting" OR "soil health assessment" OR "precision agriculture" OR "sustainable farming practices" OR "weather pattern analysis" OR (crop yield prediction) OR "crop variety selection" OR "irrigation optimization" OR "technology adoption in agriculture" OR "agroecology" OR "data-driven agriculture" OR "economic viability of crops" OR ' + "`n" + '(("predictive modeling" OR "risk assessment" OR "climate change" OR "machine learning" OR "remote sensing" OR "land use change" OR "environmental impact assessment" OR "food security" OR "agricultural policies" OR "nutrient management" OR "satellite imagery" OR "farming systems analysis" OR "earth observation data" OR "phenotyping" OR "big data analytics" OR "climate resilience") AND ((Crop Yield)))' + "`n"
$ws.Range("B10").Value = $b10
$ws.Range("C10").Value = 0.652
$ws.Range("D10").Value = 0.106
$ws.Range("E10").Value = 0.321
$ws.Range("F10").Value = 0.528
$ws.Range("G10").Value = 0.623
$ws.Rows.Item(10).AutoFit()

# --- Row 11: Baseline (Crop Yield Prediction) ---
$ws.Range("A3").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Baseline"

$ws.Range("B11").Value = "Crop Yield Prediction"
$ws.Range("C11").Value = 0.543
$ws.Range("D11").Value = 0.37
$ws.Range("E11").Value = 0.497
$ws.Range("F11").Value = 0.508
$ws.Range("G11").Value = 0.536

Write-Output "done"
